$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw counts (row 2)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 1

# Update the derived probabilities (row 5)
$ws.Range("B5").Value = 2/3
$ws.Range("D5").Value = 1/3
